$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.3357664143492041
$ws.Range("D2").Value = 0.7391102547729465

# Row 3
$ws.Range("C3").Value = -0.8729098101749544
$ws.Range("D3").Value = 0.3888356242339921

# Row 4
$ws.Range("C4").Value = -1.711101559166005
$ws.Range("D4").Value = 0.09617269127639672
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = -2.311077132584924
$ws.Range("D5").Value = 0.02702252637953406

# Row 6
$ws.Range("C6").Value = -0.2517423737536986
$ws.Range("D6").Value = 0.8027556796080899

# Row 7
$ws.Range("C7").Value = -1.2151560373487
$ws.Range("D7").Value = 0.2326778232093949

# Row 8
$ws.Range("C8").Value = -1.897962250736629
$ws.Range("D8").Value = 0.06621788816709984
$ws.Range("G8").Value = "No"

# Row 9
$ws.Range("C9").Value = -1.648157822640913
$ws.Range("D9").Value = 0.1085321923500313

# Row 10
$ws.Range("C10").Value = -2.202916773033312
$ws.Range("D10").Value = 0.03447818809248693

# Row 11
$ws.Range("C11").Value = -0.6161472385992606
$ws.Range("D11").Value = 0.5419008123366111
